# Finish ASIC control signal table: add a new "mem_out_en" control-signal
# column (S) to the "asic control signal decode" sheet, mirroring the
# formatting of the existing "mem_write_en" column (R), and nudge the
# sheet view (scroll position / selection) as recorded by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("asic control signal decode")

# --- New column header (row 1) ---------------------------------------
$ws.Cells.Item(1, 19).Value = "mem_out_en"

# --- New column data (rows 2-19) --------------------------------------
# Same pattern as column R (all zero) except row 12, which is 1.
$sValues = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $sValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 19).Value = $sValues[$i]
}

# --- Match formatting/width of column R (mem_write_en) ----------------
$ws.Range("R1:R19").Copy() | Out-Null
$ws.Range("S1:S19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Columns.Item(19).ColumnWidth = $ws.Columns.Item(18).ColumnWidth

# --- Sheet view: scroll a bit left and move the active selection ------
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$ws.Range("S13").Select() | Out-Null

$wb.Save()
